# Results from July 12, 2020 08:40:04 PM America/Chicago TZ run
# Apply the data updates described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Num($cellRef, $value) {
    $ws.Range($cellRef).Value = $value
}

function Set-Bool($cellRef, $value) {
    $ws.Range($cellRef).Value = [bool]$value
}

function Set-Text($cellRef, $value) {
    $ws.Range($cellRef).Value = $value
}

function Set-NumericText($cellRef, $value) {
    # Force a numeric-looking string to be stored as text (not a number),
    # without leaving a stray "quote prefix" style on the cell.
    $ws.Range($cellRef).Value = "'" + $value
    $ws.Range($cellRef).Style = "Normal"
}

function Clear-Cells($rangeRef) {
    $ws.Range($rangeRef).ClearContents()
    $ws.Range($rangeRef).ClearFormats()
}

# ---------------------------------------------------------------
# Row 2 - Arkansas
# ---------------------------------------------------------------
Set-Num "C2" 28367
Set-Num "D2" 321
Set-Num "E2" 6042
Set-Num "F2" 78
Set-Num "G2" 24.84
Set-Num "H2" 26.09
Set-Num "K2" 24327
Set-Num "L2" 299

# ---------------------------------------------------------------
# Row 3 - Massachusetts
# ---------------------------------------------------------------
Set-Num "B3" 44024
Set-Num "C3" 111597
Set-Num "D3" 8325
Set-Num "E3" 10493
Set-Num "G3" 9.4
Set-Num "H3" 8.22

# ---------------------------------------------------------------
# Row 5 - Iowa (now errors, data wiped)
# ---------------------------------------------------------------
Clear-Cells "B5:H5"
Set-Bool "I5" $false
Set-Bool "J5" $false
Clear-Cells "K5:L5"
Set-Text "O5" "An error occurred. ... AssertionError('No percentage deaths found.')"

# ---------------------------------------------------------------
# Row 6 - Illinois
# ---------------------------------------------------------------
Set-Num "B6" 44024
Set-Num "C6" 153916
Set-Num "D6" 7187
Set-Num "E6" 25828
Set-Num "F6" 1990
Set-Num "G6" 16.78
Set-Num "H6" 27.69

# ---------------------------------------------------------------
# Row 8 - Texas -- Bexar County
# ---------------------------------------------------------------
Set-Num "B8" 44024
Set-Num "C8" 19648
Set-Num "D8" 184

# ---------------------------------------------------------------
# Row 10 - Kentucky
# ---------------------------------------------------------------
Set-Num "B10" 44024
Set-Num "C10" 19389
Set-Num "D10" 625
Set-Num "E10" 2072
Set-Num "G10" 15.91
Set-Num "H10" 4.12
Set-Num "K10" 13029
Set-Num "L10" 587

# ---------------------------------------------------------------
# Row 12 - Wisconsin
# ---------------------------------------------------------------
Set-Num "B12" 44024
Set-Num "C12" 36448
Set-Num "D12" 820
Set-Num "E12" 6076
Set-Num "G12" 18.5
Set-Num "H12" 23.76
Set-Num "K12" 32836
Set-Num "L12" 808

# ---------------------------------------------------------------
# Row 14 - Tennessee
# ---------------------------------------------------------------
Set-Num "B14" 44024
Set-Num "C14" 61960
Set-Num "D14" 741
Set-Num "E14" 12551
Set-Num "F14" 263
Set-Num "G14" 20.26
Set-Num "H14" 35.49

# ---------------------------------------------------------------
# Row 16 - Utah (C/D/E are stored as text)
# ---------------------------------------------------------------
Set-Num "B16" 44024
Set-NumericText "C16" "29484"
Set-NumericText "D16" "215"
Set-NumericText "E16" "747"

# ---------------------------------------------------------------
# Row 18 - Missouri (now errors, data wiped)
# ---------------------------------------------------------------
Clear-Cells "B18:H18"
Set-Bool "J18" $false
Clear-Cells "K18:L18"
Set-Text "O18" "An error occurred. ... AssertionError('Unable to find ArcGIS ID 554ada3bc8b147abad21ae23d4a7ba3a')"

# ---------------------------------------------------------------
# Row 19 - California (previously errored, now succeeds)
# ---------------------------------------------------------------
Set-Num "B19" 44023
$ws.Range("B19").NumberFormat = "YYYY-MM-DD"
Set-Num "C19" 320804
Set-Num "D19" 6989
Set-Num "E19" 9021
Set-Num "F19" 613
Set-Num "G19" 4.4
Set-Num "H19" 8.9
Set-Bool "J19" $true
Set-Num "K19" 206109
Set-Num "L19" 6888
Set-Text "O19" "Success!"

# ---------------------------------------------------------------
# Row 20 - New Mexico
# ---------------------------------------------------------------
Set-Num "B20" 44024
Set-Num "C20" 15028
Set-Num "D20" 545
Set-Num "E20" 282
Set-Num "G20" 1.88

# ---------------------------------------------------------------
# Row 21 - Alaska
# ---------------------------------------------------------------
Set-Num "B21" 44024
Set-Num "C21" 1479
Set-Num "E21" 35
Set-Num "G21" 1.4
Set-Num "K21" 2499

# ---------------------------------------------------------------
# Row 22 - Texas
# ---------------------------------------------------------------
Set-Num "B22" 44024
Set-Num "C22" 25438
Set-Num "D22" 719
Set-Num "E22" 2914
Set-Num "G22" 0.11
Set-Num "H22" 0.13

# ---------------------------------------------------------------
# Row 24 - California - San Diego
# ---------------------------------------------------------------
Set-Num "B24" 44024
Set-Num "C24" 19929
Set-Num "E24" 736
Set-Num "G24" 4.65
Set-Num "K24" 15828

# ---------------------------------------------------------------
# Row 28 - California - Los Angeles
# ---------------------------------------------------------------
Set-Num "B28" 44023
Set-Num "C28" 133549
Set-Num "D28" 3809
Set-Num "E28" 3554
Set-Num "F28" 387
Set-Num "G28" 4.7
Set-Num "H28" 10.92
Set-Num "K28" 75586
Set-Num "L28" 3543

# ---------------------------------------------------------------
# Row 33
# ---------------------------------------------------------------
Set-Num "B33" 44024
Set-Num "C33" 21172
Set-Num "E33" 1252
Set-Num "G33" 7.62
Set-Num "H33" 8.06
Set-Num "K33" 16440
Set-Num "L33" 273

# ---------------------------------------------------------------
# Row 34 - Georgia
# ---------------------------------------------------------------
Set-Num "B34" 44024
Set-Num "C34" 116926
Set-Num "D34" 3001
Set-Num "E34" 31278
Set-Num "F34" 1398
Set-Num "G34" 26.75
Set-Num "H34" 46.58

# ---------------------------------------------------------------
# Row 36 - Idaho
# ---------------------------------------------------------------
Set-Num "B36" 44024
Set-Num "C36" 10902
Set-Num "E36" 150
Set-Num "G36" 1.38

# ---------------------------------------------------------------
# Row 38 - Colorado
# ---------------------------------------------------------------
Set-Num "B38" 44024
Set-Num "C38" 36913
Set-Num "E38" 1890
Set-Num "G38" 6.33
Set-Num "K38" 29866
Set-Num "L38" 1661

# ---------------------------------------------------------------
# Row 41 - Michigan
# ---------------------------------------------------------------
Set-Num "B41" 44024
Set-Num "C41" 69250
Set-Num "D41" 5984
Set-Num "E41" 20548
Set-Num "G41" 29.67
Set-Num "H41" 39.94

Write-Host "Applied covid_disparities_output update for 2020-07-12 run."
